# IFRS list fix — replace the (incorrectly scaled / stale) financial figures
# for 2014-2018 (rows 2-6) with the corrected values, and drop the forecast
# rows (2019E-2021E, rows 7-9) whose data columns were invalid, keeping
# only the "annual" / period labels (columns A-C) on those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D(4) .. AJ(36) hold the numeric series for each annual-period row.
# For each row we give one value per column in that order; `$null` means
# "leave the cell absent" (it was never populated) and '__CLEAR__' means
# "remove this cell's content" (present before, gone after).

$rowData = @{
    2 = @(2646,50,50,259,215,215,0,2963,511,2452,2446,6,82,140,207,-20,2,137,'__CLEAR__',1.88,8.130000000000001,9.16,7.54,20.83,2786.92,1164,18.33,13574,1.65,178,0.84,14.9,18466948)
    3 = @(2250,29,29,229,190,190,0,3029,413,2616,2610,6,82,100,-52,-32,10,90,'__CLEAR__',1.31,8.449999999999999,7.5,6.35,15.77,2979.18,1027,20.99,14482,1.56,245,1.14,23.22,18466948)
    4 = @(2240,52,52,296,245,245,0,3351,484,2867,2861,6,82,40,137,5,6,34,'__CLEAR__',2.33,10.95,8.949999999999999,7.68,16.89,3273.26,1326,11.88,15575,1.06,245,1.56,18.38,18466948)
    5 = @(2594,99,99,348,286,285,1,3626,530,3096,3090,7,82,13,-212,-45,4,8,0,3.82,11.01,9.58,8.19,17.11,3567.6,1544,10.24,16822,0.98,245,1.55,15.79,18466948)
    6 = @(2693,68,68,364,289,289,$null,3865,520,3345,3338,$null,82,89,-124,-57,8,81,0,2.53,10.75,8.99,7.72,15.54,3859.99,1564,8.58,18283,0.77,268,1.99,16.92,18466948)
}

foreach ($r in 2..6) {
    $values = $rowData[$r]
    $col = 4
    foreach ($v in $values) {
        if ($v -eq $null) {
            # column never populated on this row — nothing to do
        } elseif ($v -eq '__CLEAR__') {
            $ws.Cells.Item($r, $col).ClearContents()
        } else {
            $ws.Cells.Item($r, $col).Value = $v
        }
        $col = $col + 1
    }
}

# Rows 7-9 (2019E/2020E/2021E forecast rows): wipe all the numeric columns,
# keeping only A (index), B ("연간") and C (period label).
$ws.Range("D7:AJ9").ClearContents()

Write-Output "IFRS list corrected"
